# Update the module number shown at the bottom of every slide.
# This text lives in the Slide Master (it is not duplicated per-slide),
# in the shape named "TextBox 10" which currently reads "Module 5 ".

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

$target = $null
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 10") {
        $target = $sh
        break
    }
}

if ($target -ne $null) {
    $target.TextFrame.TextRange.Text = "Module 4 "
}
